# KDR003 template: increase the top page margin from 2.0cm to 2.5cm
# (0.98425196850393704 in == 70.86614173228347 pt) so the printed header
# has more clearance. This is the only user-visible change in the
# corresponding workbook revision - fix bug kdr003 #121388 #121390.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.PageSetup.TopMargin = 70.86614173228347
